$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to all cells being updated so that numeric-looking
# strings (e.g. "1.00", "7.15", "0.998") are preserved exactly as text and are
# not auto-converted to floating point numbers by Excel.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.041.20'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.527.76'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '598.54'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '134.37'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -1.54%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.526.95'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.124'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.32%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.15'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +3.61%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.125.52'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.45'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.522.55'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '65.006.52'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.07'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.26%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.42'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.95%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.69'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -1.86%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '391.40'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.576'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.669.62'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '74.27'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000114'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +1.24%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.64'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +21.80%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.76'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.92%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.93%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.530.16'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +1.47%  '
$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.58'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +2.37%  '
$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.19'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +5.50%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '168.71'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.87%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.76%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +2.92%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.822'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'ONDO'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.24'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +4.01%  '
$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '42.69'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '25.28'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -4.20%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.43'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.66'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.91'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.08%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.410.25'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.896'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +6.02%  '
